# Flora 001.xlsx - "Add unipa transcriptions (-Hasan 026)"
#
# The UNIPA sampling-code transcriptions that previously read
# "...AF001..." are corrected to "...AM001..." everywhere they occur:
#   - "Sampling Events" sheet, columns A:B (rows 2-3)
#   - "Occurrences" sheet, columns A:B (rows 2-40)
#   - "Literature" sheet, column A (row 2)
#
# A straightforward text Replace across the used columns reproduces the
# edit (and, as a side effect, the shared-string table gets rebuilt in
# the same relative order as the authoring tool produced).

$wb = $excel.ActiveWorkbook

# --- Sampling Events ---------------------------------------------------
$wsSampling = $wb.Worksheets.Item("Sampling Events")
$wsSampling.Range("A:B").Replace("AF001", "AM001") | Out-Null

# --- Occurrences --------------------------------------------------------
$wsOcc = $wb.Worksheets.Item("Occurrences")
$wsOcc.Range("A:B").Replace("AF001", "AM001") | Out-Null

# --- Literature ----------------------------------------------------------
$wsLit = $wb.Worksheets.Item("Literature")
$wsLit.Range("A:A").Replace("AF001", "AM001") | Out-Null

# --- Restore / update the view state ------------------------------------
# Sampling Events: selection widened to the two used columns.
$wsSampling.Activate()
$wsSampling.Range("A1:B1048576").Select() | Out-Null

# Occurrences: selection widened to the two used columns, still scrolled
# to row 28.
$wsOcc.Activate()
$wsOcc.Range("A1:B1048576").Select() | Out-Null

# Literature: becomes the active/selected tab, selection widened to the
# single used column.
$wsLit.Activate()
$wsLit.Range("A1:A1048576").Select() | Out-Null

Write-Host "Unipa transcriptions updated (AF001 -> AM001)."
